# Applies the "make some small changes" edit to the Journal options bullet
# list in notes/DiscussionsAndToDo.docx.

$d = $word.ActiveDocument

function ReplaceText($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# 1) Nature Scientific Reports bullet
ReplaceText "Nature Scientific Reports (where Florian’s paper 2016 was published, not discussed yet)" "Nature Scientific Reports (where Florian’s paper 2016 was published, Arup says maybe)"

# 2) PNAS bullet
ReplaceText "PNAS (suggestion from Navish, not discussed yet)" "PNAS (suggestion from Navish, Arup says it might not fly there)"

# 3) PLoS Computational Biology bullet (narrow match keeps the preceding
#    "where Mann et al. ... published" run and its proofErr markers intact)
ReplaceText ", not discussed yet)" ", Arup says maybe)"

# 4) Journal of the Royal Society Interface bullet (this also removes the
#    "_GoBack" bookmark that used to sit between "not discussed yet" and ")";
#    it gets re-created later on the new Statistical Physics bullet).
ReplaceText "(my suggestion, not discussed yet)" "(my suggestion, Arup says  maybe)"

# 5) Mathematical Biosciences bullet (narrow match preserves the "ms"
#    proofErr-wrapped run right after it)
ReplaceText "Mathematical Biosciences (if our other " "Mathematical Biosciences (Arup: if our other "

# 6) Journal of Theoretical Biology bullet
ReplaceText ", might be too field-specific/theoretical)" ", my and Navish’s opinion: might be too field-specific/theoretical)"

# 7) Add a brand-new bullet for Journal of Statistical Physics right after the
#    Journal of Theoretical Biology bullet, re-using the same list formatting,
#    and move the "_GoBack" bookmark onto it (just before the final ")").
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Journal of Theoretical Biology") {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find Journal of Theoretical Biology paragraph"
}

$insertAt = $targetPara.Range.End
$targetPara.Range.InsertParagraphAfter()

$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $insertAt) {
        $newPara = $p
        break
    }
}
if ($newPara -eq $null) {
    throw "Could not find newly inserted paragraph"
}

$newPara.Range.InsertAfter("Journal of Statistical Physics (Mehran’s suggestion, Arup says it won’t fit there)")

$bmStart = $newPara.Range.End - 2
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
